$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "pred" column (B) values to reflect the "keep highest val strategy"
$ws.Range("B6").Value = "BBB"
$ws.Range("B9").Value = "A"
$ws.Range("B10").Value = "A"
$ws.Range("B11").Value = "A"
$ws.Range("B12").Value = "A"
$ws.Range("B33").Value = "BBB"
$ws.Range("B39").Value = "B"
$ws.Range("B51").Value = "BB"
$ws.Range("B62").Value = "A"
$ws.Range("B76").Value = "BBB"
$ws.Range("B85").Value = "B"
$ws.Range("B87").Value = "BBB"
$ws.Range("B95").Value = "A"
$ws.Range("B113").Value = "A"
$ws.Range("B116").Value = "A"
$ws.Range("B118").Value = "BBB"
$ws.Range("B140").Value = "BBB"
$ws.Range("B141").Value = "BB"
$ws.Range("B143").Value = "AA"
$ws.Range("B146").Value = "A"
$ws.Range("B156").Value = "B"
$ws.Range("B157").Value = "B"
$ws.Range("B158").Value = "BB"
$ws.Range("B159").Value = "A"
$ws.Range("B161").Value = "A"
$ws.Range("B167").Value = "A"
$ws.Range("B169").Value = "BBB"
$ws.Range("B174").Value = "BB"
$ws.Range("B180").Value = "A"
$ws.Range("B182").Value = "A"
$ws.Range("B184").Value = "B"
$ws.Range("B187").Value = "AAA"
$ws.Range("B195").Value = "BBB"
$ws.Range("B212").Value = "B"
$ws.Range("B226").Value = "BBB"
$ws.Range("B229").Value = "B"
$ws.Range("B233").Value = "BBB"
$ws.Range("B239").Value = "A"
$ws.Range("B249").Value = "B"
$ws.Range("B255").Value = "BBB"
$ws.Range("B270").Value = "BB"
$ws.Range("B273").Value = "BBB"
$ws.Range("B283").Value = "BB"
$ws.Range("B286").Value = "BBB"
$ws.Range("B288").Value = "B"
$ws.Range("B290").Value = "BBB"
$ws.Range("B291").Value = "BBB"
$ws.Range("B300").Value = "A"
$ws.Range("B301").Value = "A"
$ws.Range("B302").Value = "A"
$ws.Range("B303").Value = "BB"
$ws.Range("B311").Value = "BBB"
$ws.Range("B313").Value = "B"
$ws.Range("B342").Value = "BBB"
$ws.Range("B349").Value = "B"
$ws.Range("B360").Value = "AA"
$ws.Range("B361").Value = "A"
$ws.Range("B370").Value = "BBB"
$ws.Range("B385").Value = "BB"
$ws.Range("B387").Value = "BBB"
$ws.Range("B394").Value = "BB"
